# Re-sort the sample "MatchValue" column (B) within each "PID" group (A),
# descending by MatchValue, matching the updated sort_data.xlsx fixture.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 14.392
$ws.Range("B5").Value = 12.321
$ws.Range("B6").Value = 28.134
$ws.Range("B8").Value = 12.321
